# Done with 387. First Unique Character in a String
#
# Row 12 was an (almost) empty placeholder row containing only a red "❌"
# status marker. Fill it in as a completed entry, matching the look & feel
# of the other completed rows (e.g. row 7) by copying that row's
# formatting first, then overwriting the row's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pull the cell formatting (fill colours, wrap text, font colour, etc.)
# from the previous fully-filled-in row (row 7) onto row 12 so the new
# entry gets the same per-column styling used elsewhere in the sheet.
$ws.Range("A7:I7").Copy()
$ws.Range("A12:I12").PasteSpecial(-4122)

# Now fill in the actual problem data for #387 "First Unique Character in
# a String".
$ws.Range("A12").Value = 387
$ws.Range("B12").Value = "First Unique Character in a String"
$ws.Range("C12").Value = "String"
$ws.Range("D12").Value = "String, Hash table, Queue"
$ws.Range("E12").Value = "Dict"
$ws.Range("F12").Value = "Easy"
$ws.Range("G12").Value = 2
$ws.Range("H12").Value = "✅"
$ws.Range("I12").Value = "Given 2 solutions. 1 O(n^2) and 2nd O(n). Havent seen solution pane"

# Reflect the author's new cursor position / scroll offset in the sheet
# view (was scrolled to top with B15 selected; now scrolled down with
# I13 selected).
$excel.ActiveWindow.ScrollRow = 6
$ws.Range("I13").Select()
